# Weekly price list: a new weekly reading was recorded, which inserts a
# new row at row 60 and pushes every subsequent reading down by one row
# (the last existing row is duplicated down to the new last row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 60 - Excel shifts row 60..256
# down to 61..257 and carries formatting (incl. the date number format
# on column D) along with it.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Range("A60").Value = 3
$ws.Range("B60").Value = "Femacal de La Calera"
$ws.Range("C60").Value = "Coquimbo"
$ws.Range("D60").Value = 44592
$ws.Range("E60").Value = 5
$ws.Range("F60").Value = 100112039
$ws.Range("G60").Value = "Ciboulette"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 160
$ws.Range("K60").Value = 1500
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = 1500
$ws.Range("N60").Value = "`$/docena de atados"
$ws.Range("O60").Value = "Provincia de Quillota"
$ws.Range("P60").Value = 500
$ws.Range("Q60").Value = 3
$ws.Range("R60").Value = "Hortaliza"
